$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 2865.2942
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 3439.2307
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 10317.6921
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -10857.6921

# Row 73
$ws.Range("H73").Value = 2865.2942
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 3439.2307
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 10317.6921
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -12189.6921

# Row 129
$ws.Range("H129").Value = 2779072.2
$ws.Range("J129").Value = 1371.8462
$ws.Range("L129").Value = 4115.5386
$ws.Range("N129").Value = -14115.5386

# Row 137
$ws.Range("H137").Value = 6672473
$ws.Range("I137").Value = 9097372
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 27292116
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -27289566
$ws.Range("N137").Value = -17100

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 250.4
$ws.Range("J4").Value = 240.66667
$ws.Range("L4").Value = 240.66667
$ws.Range("N4").Value = -472.66667

# Row 24
$ws.Range("H24").Value = 32000
$ws.Range("J24").Value = 32000
$ws.Range("L24").Value = 32000
$ws.Range("N24").Value = -32748

# Row 36
$ws.Range("H36").Value = 38944.25
$ws.Range("I36").Value = 5366.75
$ws.Range("J36").Value = 72521.75
$ws.Range("K36").Value = 5366.75
$ws.Range("L36").Value = 72521.75
$ws.Range("M36").Value = -5020.75
$ws.Range("N36").Value = -73213.75

# Row 100
$ws.Range("H100").Value = 32000
$ws.Range("J100").Value = 32000
$ws.Range("L100").Value = 32000
$ws.Range("N100").Value = -34164

$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 20274
$ws.Range("I26").Value = 20274
$ws.Range("K26").Value = 20274
$ws.Range("M26").Value = -19982

# Row 100
$ws.Range("H100").Value = 18728.666
$ws.Range("J100").Value = 18728.666
$ws.Range("L100").Value = 18728.666
$ws.Range("N100").Value = -20892.666

$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()

# Row 21
$ws.Range("H21").Value = 46981
$ws.Range("I21").Value = 913
$ws.Range("J21").Value = 70015
$ws.Range("K21").Value = 913
$ws.Range("L21").Value = 70015
$ws.Range("M21").Value = -678
$ws.Range("N21").Value = -70485

# Row 31
$ws.Range("H31").Value = 3128529
$ws.Range("I31").Value = 3574176
$ws.Range("J31").Value = 9000
$ws.Range("K31").Value = 3574176
$ws.Range("L31").Value = 9000
$ws.Range("M31").Value = -3573881
$ws.Range("N31").Value = -9590

# Row 34
$ws.Range("H34").Value = 3128529
$ws.Range("I34").Value = 3574176
$ws.Range("J34").Value = 9000
$ws.Range("K34").Value = 3574176
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -3573974
$ws.Range("N34").Value = -9404

# Row 74
$ws.Range("H74").Value = 23677.428
$ws.Range("J74").Value = 23677.428
$ws.Range("L74").Value = 23677.428
$ws.Range("N74").Value = -25425.428

# Row 77
$ws.Range("H77").Value = 23677.428
$ws.Range("J77").Value = 23677.428
$ws.Range("L77").Value = 71032.284
$ws.Range("N77").Value = -79768.284

# Row 88
$ws.Range("H88").Value = 32150
$ws.Range("J88").Value = 32150
$ws.Range("L88").Value = 32150
$ws.Range("N88").Value = -32962

# Row 91
$ws.Range("H91").Value = 32150
$ws.Range("J91").Value = 32150
$ws.Range("L91").Value = 32150
$ws.Range("N91").Value = -34958

# Row 92
$ws.Range("H92").Value = 16500
$ws.Range("J92").Value = 16500
$ws.Range("L92").Value = 16500
$ws.Range("N92").Value = -21492

# Row 96
$ws.Range("H96").Value = 23640
$ws.Range("J96").Value = 23640
$ws.Range("L96").Value = 23640
$ws.Range("N96").Value = -29132

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1166.6129
$ws.Range("I131").Value = 4500
$ws.Range("J131").Value = 997.11865
$ws.Range("K131").Value = 13500
$ws.Range("L131").Value = 2991.35595
$ws.Range("M131").Value = -8460
$ws.Range("N131").Value = -13071.35595

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 6005016.5
$ws.Range("J11").Value = 6010042
$ws.Range("L11").Value = 6010042
$ws.Range("N11").Value = -6010320

# Row 21
$ws.Range("H21").Value = 45003.5
$ws.Range("J21").Value = 45003.5
$ws.Range("L21").Value = 45003.5
$ws.Range("N21").Value = -45349.5

# Row 30
$ws.Range("H30").Value = 45003.5
$ws.Range("J30").Value = 45003.5
$ws.Range("L30").Value = 45003.5
$ws.Range("N30").Value = -45213.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 125001470
$ws.Range("I22").Value = 250000200
$ws.Range("J22").Value = 2737.5
$ws.Range("K22").Value = 250000200
$ws.Range("L22").Value = 2737.5
$ws.Range("M22").Value = -249999905
$ws.Range("N22").Value = -3327.5

# Row 23
$ws.Range("H23").Value = 145243.72
$ws.Range("I23").Value = 168451
$ws.Range("K23").Value = 168451
$ws.Range("M23").Value = -168221

# Row 27
$ws.Range("H27").Value = 125001470
$ws.Range("I27").Value = 250000200
$ws.Range("J27").Value = 2737.5
$ws.Range("K27").Value = 250000200
$ws.Range("L27").Value = 2737.5
$ws.Range("M27").Value = -250000093
$ws.Range("N27").Value = -2951.5

$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 16625
$ws.Range("J45").Value = 16625
$ws.Range("L45").Value = 16625
$ws.Range("N45").Value = -17607

# Row 62
$ws.Range("H62").Value = 3650
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -1676
$ws.Range("N62").Value = -6248

# Row 65
$ws.Range("H65").Value = 3650
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -8380
$ws.Range("N65").Value = -31240

# Row 69
$ws.Range("H69").Value = 29625
$ws.Range("J69").Value = 29625
$ws.Range("L69").Value = 29625
$ws.Range("N69").Value = -31123

# Row 72
$ws.Range("H72").Value = 29625
$ws.Range("J72").Value = 29625
$ws.Range("L72").Value = 88875
$ws.Range("N72").Value = -96363

# Row 92
$ws.Range("H92").Value = 32449.75
$ws.Range("J92").Value = 32449.75
$ws.Range("L92").Value = 32449.75
$ws.Range("N92").Value = -37441.75

# Row 113
$ws.Range("H113").Value = 1246.4706
$ws.Range("I113").Value = 576.7778
$ws.Range("J113").Value = 1999.875
$ws.Range("K113").Value = 1730.3334
$ws.Range("L113").Value = 5999.625
$ws.Range("M113").Value = 439.6666
$ws.Range("N113").Value = -10339.625
